$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud"
$ws.Range("G3").Value = "Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G4").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G5").Value = "Dr. Hanan Ragab, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Nesma"
$ws.Range("G6").Value = "Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud"
$ws.Range("G7").Value = "Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G8").Value = "Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Administrator, Dr. Manar Montaser, Dr. Asmaa Reda"
$ws.Range("G9").Value = "Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G10").Value = "Dr. Shimaa Ahmad Mekki, Dr. Sara Wael, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali"
$ws.Range("G11").Value = "Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G13").Value = "Dr. Safa Hany, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Range("G15").Value = "Dr. Amal Awwad, D Wessam Atef"
$ws.Range("G16").Value = "Dr. Amal Awwad, Dr. Nourhan Mohammad"
$ws.Range("G17").Value = "Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen"
$ws.Range("G23").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G24").Value = "Dr. Youstina Magdy, Dr. Maryam Ashraf, Dr. Aya Emad, Dr. Monica, Dr. Salma Hassan, Dr. Marina Atef, Dr. Remon, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida"
$ws.Range("G25").Value = "Dr. Youstina Magdy, Dr. Aya Emad, Dr. Marina Atef, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah"
$ws.Range("G27").Value = "Dr. Eman Mohammad Al, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida"
$ws.Range("G28").Value = "Dr. Aya Hanafy, Dr. Nardine, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Neveen Nashaat, Dr. Wafaa Ebida"
$ws.Range("G29").Value = "Dr. Monica, Dr. Eman Samir Gabry, Dr. Remon, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
$ws.Range("G30").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G31").Value = "Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G32").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G33").Value = "Dr. Hanan Ragab, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Nesma"
$ws.Range("G34").Value = "Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud"
$ws.Range("G35").Value = "Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G36").Value = "Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Administrator, Dr. Manar Montaser, Dr. Asmaa Reda"
$ws.Range("G37").Value = "Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud"
$ws.Range("G38").Value = "Dr. Shimaa Ahmad Mekki, Dr. Sara Wael, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali"
$ws.Range("G39").Value = "Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G41").Value = "Dr. Safa Hany, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Range("G43").Value = "Dr. Amal Awwad, D Wessam Atef"
$ws.Range("G44").Value = "Dr. Amal Awwad, Dr. Nourhan Mohammad"
$ws.Range("G45").Value = "Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Sarah Abdelmohsen"
$ws.Range("G51").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G52").Value = "Dr. Youstina Magdy, Dr. Maryam Ashraf, Dr. Aya Emad, Dr. Monica, Dr. Salma Hassan, Dr. Marina Atef, Dr. Remon, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida"
$ws.Range("G53").Value = "Dr. Youstina Magdy, Dr. Aya Emad, Dr. Marina Atef, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah"
$ws.Range("G55").Value = "Dr. Eman Mohammad Al, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida"
$ws.Range("G56").Value = "Dr. Aya Hanafy, Dr. Nardine, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody, Dr. Neveen Nashaat, Dr. Wafaa Ebida"
$ws.Range("G57").Value = "Dr. Monica, Dr. Eman Samir Gabry, Dr. Remon, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
